$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44784
$ws.Range("J2").Value = 520
$ws.Range("K2").Value = 11500
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11750
$ws.Range("P2").Value = 294

$ws.Range("D3").Value = 44484
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9500
$ws.Range("O3").Value = 'Provincia del Elquí'
$ws.Range("P3").Value = 317

$ws.Range("D4").Value = 44858
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 9500
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 9750
$ws.Range("P4").Value = 325

$ws.Range("D5").Value = 45070
$ws.Range("H5").Value = 'Madrigal'
$ws.Range("J5").Value = 360
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17500
$ws.Range("N5").Value = '$/caja 40 unidades'
$ws.Range("O5").Value = 'Provincia del Elquí'
$ws.Range("P5").Value = 438
$ws.Range("Q5").Value = 40

$ws.Range("D6").Value = 45037
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 16500
$ws.Range("P6").Value = 412

$ws.Range("D7").Value = 45069
$ws.Range("J7").Value = 500

$ws.Range("D8").Value = 44687
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 19000
$ws.Range("M8").Value = 18500
$ws.Range("P8").Value = 617

$ws.Range("D9").Value = 45049
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("O9").Value = 'Provincia del Elquí'
$ws.Range("P9").Value = 438

$ws.Range("D10").Value = 44839
$ws.Range("H10").Value = 'Española'
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 12500
$ws.Range("N10").Value = '$/caja 30 unidades'
$ws.Range("P10").Value = 417
$ws.Range("Q10").Value = 30

$ws.Range("D11").Value = 44767
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 362

$ws.Range("D12").Value = 44701
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 19000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 19500
$ws.Range("P12").Value = 650

$ws.Range("D13").Value = 44420
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14500
$ws.Range("O13").Value = 'Provincia de Limarí'
$ws.Range("P13").Value = 362

$ws.Range("D14").Value = 44420
$ws.Range("J14").Value = 700
$ws.Range("K14").Value = 13000
$ws.Range("M14").Value = 13500
$ws.Range("O14").Value = 'Provincia del Elquí'
$ws.Range("P14").Value = 338

$ws.Range("D15").Value = 44729
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 16000
$ws.Range("M15").Value = 16500
$ws.Range("P15").Value = 412

$ws.Range("D16").Value = 44855
$ws.Range("H16").Value = 'Española'
$ws.Range("J16").Value = 540
$ws.Range("K16").Value = 9500
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 9750
$ws.Range("N16").Value = '$/caja 30 unidades'
$ws.Range("O16").Value = 'Provincia del Elquí'
$ws.Range("P16").Value = 325
$ws.Range("Q16").Value = 30

$ws.Range("D17").Value = 45082
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("P17").Value = 438

$ws.Range("D18").Value = 44438
$ws.Range("K18").Value = 11000
$ws.Range("M18").Value = 11500
$ws.Range("P18").Value = 383

$ws.Range("D19").Value = 44790
$ws.Range("H19").Value = 'Española'
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14500
$ws.Range("N19").Value = '$/caja 30 unidades'
$ws.Range("O19").Value = 'Provincia de Limarí'
$ws.Range("P19").Value = 483
$ws.Range("Q19").Value = 30

$ws.Range("D20").Value = 44790
$ws.Range("K20").Value = 11500
$ws.Range("M20").Value = 11750
$ws.Range("O20").Value = 'Provincia del Elquí'
$ws.Range("P20").Value = 294

$ws.Range("D21").Value = 45079
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 16500
$ws.Range("M21").Value = 16750
$ws.Range("P21").Value = 419

$ws.Range("D22").Value = 44498
$ws.Range("K22").Value = 8500
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 8750
$ws.Range("P22").Value = 292

$ws.Range("D23").Value = 44427
$ws.Range("H23").Value = 'Madrigal'
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 12500
$ws.Range("N23").Value = '$/caja 40 unidades'
$ws.Range("O23").Value = 'Provincia de Limarí'
$ws.Range("P23").Value = 312
$ws.Range("Q23").Value = 40

$ws.Range("D24").Value = 45090
$ws.Range("H24").Value = 'Madrigal'
$ws.Range("J24").Value = 340
$ws.Range("K24").Value = 15500
$ws.Range("L24").Value = 16000
$ws.Range("M24").Value = 15750
$ws.Range("N24").Value = '$/caja 40 unidades'
$ws.Range("P24").Value = 394
$ws.Range("Q24").Value = 40

$ws.Range("D25").Value = 44426
$ws.Range("H25").Value = 'Española'
$ws.Range("J25").Value = 600
$ws.Range("K25").Value = 11500
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = 11750
$ws.Range("N25").Value = '$/caja 30 unidades'
$ws.Range("O25").Value = 'Provincia de Limarí'
$ws.Range("P25").Value = 392
$ws.Range("Q25").Value = 30

$ws.Range("D26").Value = 44426
$ws.Range("H26").Value = 'Madrigal'
$ws.Range("K26").Value = 12500
$ws.Range("L26").Value = 13000
$ws.Range("M26").Value = 12750
$ws.Range("N26").Value = '$/caja 40 unidades'
$ws.Range("O26").Value = 'Provincia de Limarí'
$ws.Range("P26").Value = 319
$ws.Range("Q26").Value = 40
